# Auto-generated Excel COM-interop script applying the Behemoth_Profits diff.
# For each changed cell we set .Value directly; cells that become empty (removed
# from the XML in the diff) are cleared with ClearContents so the <c> element drops
# out of the saved worksheet, matching the target OOXML exactly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1079.4286
$ws.Range("I41").Value = 739.63635
$ws.Range("K41").Value = 739.63635
$ws.Range("M41").Value = -299.63635
$ws.Range("H43").Value = 2752.55
$ws.Range("I43").Value = 1318.4736
$ws.Range("K43").Value = 1318.4736
$ws.Range("M43").Value = -1249.4736
$ws.Range("H64").Value = 4676.4707
$ws.Range("I64").Value = 3625
$ws.Range("K64").Value = 3625
$ws.Range("M64").Value = -3377
$ws.Range("H67").Value = 4676.4707
$ws.Range("I67").Value = 3625
$ws.Range("K67").Value = 3625
$ws.Range("M67").Value = -2767
$ws.Range("H74").Value = 17821.143
$ws.Range("I74").Value = 23389.6
$ws.Range("J74").Value = 3900
$ws.Range("K74").Value = 23389.6
$ws.Range("L74").Value = 3900
$ws.Range("M74").Value = -22453.6
$ws.Range("N74").Value = -5772
$ws.Range("H77").Value = 17821.143
$ws.Range("I77").Value = 23389.6
$ws.Range("J77").Value = 3900
$ws.Range("K77").Value = 116948
$ws.Range("L77").Value = 19500
$ws.Range("M77").Value = -112268
$ws.Range("N77").Value = -28860
$ws.Range("H86").Value = 7446.8237
$ws.Range("I86").Value = 6673.067
$ws.Range("J86").Value = 13250
$ws.Range("K86").Value = 6673.067
$ws.Range("L86").Value = 13250
$ws.Range("M86").Value = -5550.067
$ws.Range("N86").Value = -15496
$ws.Range("H89").Value = 7446.8237
$ws.Range("I89").Value = 6673.067
$ws.Range("J89").Value = 13250
$ws.Range("K89").Value = 33365.335
$ws.Range("L89").Value = 66250
$ws.Range("M89").Value = -27749.335
$ws.Range("N89").Value = -77482
$ws.Range("H96").Value = 1480.875
$ws.Range("I96").Value = 1411.6
$ws.Range("J96").Value = 1596.3334
$ws.Range("K96").Value = 4234.799999999999
$ws.Range("L96").Value = 4789.0002
$ws.Range("M96").Value = -2861.799999999999
$ws.Range("N96").Value = -7535.0002
$ws.Range("H98").Value = 45455812
$ws.Range("I98").Value = 45455812
$ws.Range("K98").Value = 45455812
$ws.Range("M98").Value = -45454314
$ws.Range("H122").Value = 45455812
$ws.Range("I122").Value = 45455812
$ws.Range("K122").Value = 136367436
$ws.Range("M122").Value = -136364986
$ws.Range("H132").Value = 1859.762
$ws.Range("I132").Value = 1950.3158
$ws.Range("J132").Value = 999.5
$ws.Range("K132").Value = 5850.9474
$ws.Range("L132").Value = 2998.5
$ws.Range("M132").Value = -3320.9474
$ws.Range("N132").Value = -8058.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 2796.6667
$ws.Range("I30").Value = 695
$ws.Range("K30").Value = 695
$ws.Range("M30").Value = -545
$ws.Range("H45").Value = 2234.7058
$ws.Range("I45").Value = 1750
$ws.Range("K45").Value = 1750
$ws.Range("M45").Value = -1373
$ws.Range("H61").Value = 23861630
$ws.Range("I61").Value = 41668344
$ws.Range("K61").Value = 41668344
$ws.Range("M61").Value = -41668132
$ws.Range("H97").Value = 589.8333
$ws.Range("I97").Value = 721.5714
$ws.Range("J97").Value = 128.75
$ws.Range("K97").Value = 721.5714
$ws.Range("L97").Value = 128.75
$ws.Range("M97").Value = -225.5714
$ws.Range("N97").Value = -1120.75
$ws.Range("H102").Value = 3071.3076
$ws.Range("J102").Value = 3999.5
$ws.Range("L102").Value = 3999.5
$ws.Range("N102").Value = -7243.5
$ws.Range("H132").Value = 4772.4
$ws.Range("J132").Value = 8219.666999999999
$ws.Range("L132").Value = 24659.001
$ws.Range("N132").Value = -29719.001
$ws.Range("H136").Value = 23861630
$ws.Range("I136").Value = 41668344
$ws.Range("K136").Value = 125005032
$ws.Range("M136").Value = -125002482

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1406.4
$ws.Range("I94").Value = 1461.3077
$ws.Range("J94").Value = 1049.5
$ws.Range("K94").Value = 1461.3077
$ws.Range("L94").Value = 1049.5
$ws.Range("M94").Value = -1010.3077
$ws.Range("N94").Value = -1951.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 16150.333
$ws.Range("I6").Value = 6250
$ws.Range("K6").Value = 6250
$ws.Range("M6").Value = -6137
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H25").Value = 2010
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H105").Value = 1982
$ws.Range("I105").Value = 1565.6
$ws.Range("K105").Value = 1565.6
$ws.Range("M105").Value = 181.4000000000001
$ws.Range("H120").Value = 28500
$ws.Range("J120").Value = 28500
$ws.Range("L120").Value = 28500
$ws.Range("N120").Value = -35758
$ws.Range("H132").Value = 2867.75
$ws.Range("I132").Value = 2941.3
$ws.Range("K132").Value = 8823.900000000001
$ws.Range("M132").Value = -6293.900000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1552.75
$ws.Range("I122").Value = 793.4
$ws.Range("J122").Value = 1897.909
$ws.Range("K122").Value = 7140.599999999999
$ws.Range("L122").Value = 17081.181
$ws.Range("M122").Value = -4690.599999999999
$ws.Range("N122").Value = -21981.181
$ws.Range("H141").Value = 10984.214
$ws.Range("I141").Value = 9847.375
$ws.Range("K141").Value = 29542.125
$ws.Range("M141").Value = -24362.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10499301
$ws.Range("I11").Value = 17307500
$ws.Range("J11").Value = 4542126
$ws.Range("K11").Value = 17307500
$ws.Range("L11").Value = 4542126
$ws.Range("M11").Value = -17307361
$ws.Range("N11").Value = -4542404
$ws.Range("H18").Value = 30000
$ws.Range("I18").Value = 30000
$ws.Range("K18").Value = 30000
$ws.Range("M18").Value = -29707
$ws.Range("H21").Value = 19499.5
$ws.Range("I21").Value = 19499.5
$ws.Range("K21").Value = 19499.5
$ws.Range("M21").Value = -19326.5
$ws.Range("H30").Value = 19499.5
$ws.Range("I30").Value = 19499.5
$ws.Range("K30").Value = 19499.5
$ws.Range("M30").Value = -19394.5
$ws.Range("H97").Value = 679.25
$ws.Range("I97").Value = 722
$ws.Range("J97").Value = 38
$ws.Range("K97").Value = 722
$ws.Range("L97").Value = 38
$ws.Range("M97").Value = -226
$ws.Range("N97").Value = -1030
$ws.Range("H126").Value = 1970
$ws.Range("I126").Value = 1975
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 5925
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -3455
$ws.Range("N126").Value = -10790
$ws.Range("H132").Value = 125001680
$ws.Range("I132").Value = 166668400
$ws.Range("K132").Value = 500005200
$ws.Range("M132").Value = -500002670
$ws.Range("H136").Value = 16949.334
$ws.Range("J136").Value = 16949.334
$ws.Range("L136").Value = 50848.00199999999
$ws.Range("N136").Value = -55948.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H20").Value = 33268.668
$ws.Range("I20").Value = 9800
$ws.Range("J20").Value = 45003
$ws.Range("K20").Value = 9800
$ws.Range("L20").Value = 45003
$ws.Range("M20").Value = -9574
$ws.Range("N20").Value = -45455
$ws.Range("H43").Value = 2537747.5
$ws.Range("I43").Value = 3372000
$ws.Range("J43").Value = 34989.5
$ws.Range("K43").Value = 3372000
$ws.Range("L43").Value = 34989.5
$ws.Range("M43").Value = -3371807
$ws.Range("N43").Value = -35375.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1491.6923
$ws.Range("J96").Value = 550
$ws.Range("L96").Value = 550
$ws.Range("N96").Value = -3296
$ws.Range("H122").Value = 8778.25
$ws.Range("I122").Value = 6878
$ws.Range("J122").Value = 10135.571
$ws.Range("K122").Value = 20634
$ws.Range("L122").Value = 30406.713
$ws.Range("M122").Value = -18184
$ws.Range("N122").Value = -35306.713
$ws.Range("H126").Value = 1464.091
$ws.Range("I126").Value = 1472.7778
$ws.Range("K126").Value = 4418.3334
$ws.Range("M126").Value = -1948.3334

Write-Host "Applied 210 cell updates and 3 clears."
